$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A (shifts old A:N -> B:O, including the "Компания" header
# that was in A1 moving to B1, along with column A's old custom width).
$ws.Columns("A").Insert()

# Copy formatting (font/border/alignment) from the shifted header cell (now B1, which
# carries the original A1 header style) back onto the freshly inserted A1 before we
# restore its text, so A1 keeps its original bordered/centered header style.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# Restore the "Компания" header back into A1, and put the new "PES" header into B1.
$ws.Range("A1").Value = "Компания"
$ws.Range("B1").Value = "PES"

# Match the final selection state left behind in the sheet.
$ws.Range("B2").Select()
